$d = $word.ActiveDocument

$replacements = @(
    @("55÷6=", "66÷8="),
    @("79÷9=", "60÷9="),
    @("96÷9=", "72÷8="),
    @("97÷9=", "29÷9="),
    @("66÷5=", "74÷6="),
    @("28÷9=", "89÷8="),
    @("73÷8=", "91÷8="),
    @("78÷6=", "77÷5="),
    @("59÷5=", "53÷3="),
    @("62÷9=", "65÷3="),
    @("18÷3=", "82÷4="),
    @("95÷3=", "90÷8="),
    @("22÷4=", "59÷7="),
    @("20÷7=", "39÷3="),
    @("58÷3=", "79÷5="),
    @("95÷9=", "34÷8="),
    @("40÷7=", "73÷4="),
    @("35÷7=", "11÷3="),
    @("61÷8=", "77÷9="),
    @("47÷4=", "32÷9="),
    @("37÷8=", "80÷5="),
    @("35÷2=", "54÷6="),
    @("76÷9=", "27÷3="),
    @("54÷2=", "28÷4="),
    @("29÷7=", "83÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
